# Update the "Förändrad" (Changed) date column C for rows 2-150
# from serial date 45189 (2023-09-20) to 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 150; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
